$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows (top-down, using row numbers valid at each step)
$ws.Rows(3).Insert()
$ws.Rows(6).Insert()
$ws.Rows(9).Insert()
$ws.Rows(11).Insert()
$ws.Rows(12).Insert()

# Fill in final values for rows 2-15 (columns A-J)
# Row 2
$ws.Range("A2").Value = "education"
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = "PhD. Candidate in Dragon Riding"
$ws.Range("D2").Value = "Basgiath War College"
$ws.Range("E2").Value = "Navarre"
$ws.Range("F2").Value = 2015
$ws.Range("G2").Value = 2020
$ws.Range("H2").Value = "Riding and learning combat while on their dragon; able to acess their powers through dragon connection"
$ws.Range("I2").Value = "NA"
$ws.Range("J2").Value = "NA"

# Row 3
$ws.Range("A3").Value = "education"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "Certificate For Being the Most Hated Second Year Rider"
$ws.Range("D3").Value = "Basgiath War College"
$ws.Range("E3").Value = "Navarre"
$ws.Range("F3").Value = 2016
$ws.Range("G3").Value = 2017
$ws.Range("H3").Value = "Awarded for being the most hated rider, including mulyiple attempts on her life by other unbonded riders/cadets"
$ws.Range("I3").Value = "NA"
$ws.Range("J3").Value = $null

# Row 4
$ws.Range("A4").Value = "education"
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = "B.S., Weapons and Potions Master"
$ws.Range("D4").Value = "Basgiath War College"
$ws.Range("E4").Value = "Navarre"
$ws.Range("F4").Value = 2011
$ws.Range("G4").Value = 2015
$ws.Range("H4").Value = "Honors Thesis: Defeating the most dangerous cadets by poisoning them as well as defeating in combat"
$ws.Range("I4").Value = "NA"
$ws.Range("J4").Value = "NA"

# Row 5
$ws.Range("A5").Value = "war_college_positions"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = "Second Squad, Flame Section, Fourth Wing"
$ws.Range("D5").Value = "Basgiath War College"
$ws.Range("E5").Value = "Navarre"
$ws.Range("F5").Value = 2011
$ws.Range("G5").Value = 2015
$ws.Range("H5").Value = "Bacame part of the first year Second Squad Fourth Wing crew"
$ws.Range("I5").Value = "NA"
$ws.Range("J5").Value = "NA"

# Row 6
$ws.Range("A6").Value = "war_college_positions"
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = "Voted Second Squads Most Improved Cadet"
$ws.Range("D6").Value = "Basgiath War College"
$ws.Range("E6").Value = "Navarre"
$ws.Range("F6").Value = 2011
$ws.Range("G6").Value = 2012
$ws.Range("H6").Value = "Position given as recipient was expected to fail during the first challenge and instead, became the most powerful dragon rider"
$ws.Range("I6").Value = "NA"
$ws.Range("J6").Value = "NA"

# Row 7
$ws.Range("A7").Value = "war_college_positions"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = "Iron Squad"
$ws.Range("D7").Value = "Basgiath War College"
$ws.Range("E7").Value = "Navarre"
$ws.Range("F7").Value = 2011
$ws.Range("G7").Value = 2015
$ws.Range("H7").Value = "Is part of the Riders Quadrant and is a part of the Iron Squad"
$ws.Range("I7").Value = "NA"
$ws.Range("J7").Value = "NA"

# Row 8
$ws.Range("A8").Value = "data_science_writings"
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = "How to stay one step ahead of your aopponents: A Guide"
$ws.Range("D8").Value = "Sorrengail Vlog"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = 2014
$ws.Range("G8").Value = "NA"
$ws.Range("H8").Value = "A step by step guide of how to poison your opponents before facing them on the sparring mats, winning before the match starts"
$ws.Range("I8").Value = "Restricted files"
$ws.Range("J8").Value = "NA"

# Row 9
$ws.Range("A9").Value = "data_science_writings"
$ws.Range("B9").Value = $true
$ws.Range("C9").Value = "Being Xaden Riorson's Partner: Blessing or a Cruse"
$ws.Range("D9").Value = "Sorrengail Vlog"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = 2015
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "What it is like being Xaden Riorson's partner as well as sharing mated dragons. Advise for future riders and who *not* to chose "
$ws.Range("I9").Value = "Vlog written after Xaden's post orders"
$ws.Range("J9").Value = "NA"

# Row 10
$ws.Range("A10").Value = "data_science_writings"
$ws.Range("B10").Value = $true
$ws.Range("C10").Value = "Riding a Dragon on a saddle 101"
$ws.Range("D10").Value = "Sorrengail Vlog"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = 2014
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "Tutorial of how to make a one of a kind sadle for riders who would otherwise not be able to stay on their dragon"
$ws.Range("I10").Value = "Provided advice for dealing with moody dragon"
$ws.Range("J10").Value = "NA"

# Row 11
$ws.Range("A11").Value = "about_me_press"
$ws.Range("B11").Value = $true
$ws.Range("C11").Value = "First rider to be chosen by two dragons"
$ws.Range("D11").Value = "Riders Encyclopedia"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = 2017
$ws.Range("G11").Value = "NA"
$ws.Range("H11").Value = "Story of how Violet was chosen by Tairn and Andarna, both incredibly rare dragons, not to mention Tairn's pairing with Sgaeyl"
$ws.Range("I11").Value = "NA"
$ws.Range("J11").Value = "NA"

# Row 12
$ws.Range("A12").Value = "about_me_press"
$ws.Range("B12").Value = $true
$ws.Range("C12").Value = "Sorrengail Family Tree"
$ws.Range("D12").Value = "Riders Encyclopedia"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = 2016
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "Addendum provided by scribes: Family tree was ripped out in 2017 after Violet and Mira's treason to Navarre"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"

# Row 13
$ws.Range("A13").Value = "about_me_press"
$ws.Range("B13").Value = $true
$ws.Range("C13").Value = "Daughter to General Lillith Sorrengail, sister to Mira and now deceased Brennan"
$ws.Range("D13").Value = "Riders Encyclopedia"
$ws.Range("E13").Value = "N/A"
$ws.Range("F13").Value = 2015
$ws.Range("G13").Value = "NA"
$ws.Range("H13").Value = "Story of every rider to come out of Basgiath War College Vol. 52"
$ws.Range("I13").Value = "NA"
$ws.Range("J13").Value = "NA"

# Row 14
$ws.Range("A14").Value = "by_me_press"
$ws.Range("B14").Value = $true
$ws.Range("C14").Value = "How to handle two of the most powerful dragons of all time"
$ws.Range("D14").Value = "Brennon's Diary; Violet's Addendum"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = 2016
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "First person excerpts of dealing with Tairn and Andarna provided by Violet"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"

# Row 15
$ws.Range("A15").Value = "by_me_press"
$ws.Range("B15").Value = $true
$ws.Range("C15").Value = "How to cope with having your once deceased brother be not deceased"
$ws.Range("D15").Value = "Brennon's Diary; Violet's Addendum"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = 2017
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = "First person excepts of Violet finiding out that her brother who was once though dead is actually alive and part of the rebellion you were once taught to hate and fight against"
$ws.Range("I15").Value = "NA"
$ws.Range("J15").Value = "NA"

# Re-apply the AutoFilter over the new data extent
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A5:J14").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new AutoFilter range
$n = $wb.Names.Item("positions!_FilterDatabase")
$n.RefersTo = "=positions!`$A`$5:`$J`$14"

# Restore the last active selection
$ws.Range("E31").Select()
